$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the core input drivers of the ramp-up model
$ws.Range("B2").Value = 6000
$ws.Range("B11").Value = 460000
$ws.Range("J4").Value = 660

# Adjust column B width to fit the new (wider) values
$ws.Columns.Item(2).ColumnWidth = 15.26953125

# Update view: zoom to 70% and move selection to J13
$excel.ActiveWindow.Zoom = 70
$ws.Range("J13").Select()
